{"js": "// The canonical OOXML diff for this change consists entirely of XML\n// attribute re-orderings (e.g. `w:tab w:val=\"..\" w:pos=\"..\"` becoming\n// `w:tab w:pos=\"..\" w:val=\"..\"`, namespace declarations being sorted\n// alphabetically, etc.) introduced by the tooling that produced the\n// diff. Every element's attribute *set* (names + values) is identical\n// before and after - only the serialized attribute order differs.\n// There is no visible or semantic change to the document's text,\n// formatting, structure, tab stops, page margins, footnotes, or\n// styles for Office.js (or any document-editing API) to reproduce.\n//\n// So this script intentionally performs no content mutation - it just\n// touches the body so the host can confirm the document loads/round-\n// trips cleanly, matching the (semantically empty) target state.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The canonical OOXML diff for this change consists entirely of XML\n# attribute re-orderings (e.g. `w:tab w:val=\"..\" w:pos=\"..\"` becoming\n# `w:tab w:pos=\"..\" w:val=\"..\"`, namespace declarations being sorted\n# alphabetically, etc.) introduced by the tooling that produced the\n# diff. Every element's attribute *set* (names + values) is identical\n# before and after - only the serialized attribute order differs.\n# There is no visible or semantic change to the document's text,\n# formatting, structure, tab stops, page margins, footnotes, or\n# styles for the Word object model (or any document-editing API) to\n# reproduce.\n#\n# So this script intentionally performs no content mutation - it just\n# touches the document to confirm it round-trips cleanly, matching the\n# (semantically empty) target state.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
